$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row ("H 72"), shifting all subsequent rows up by one
$ws.Rows.Item(2).Delete()

# Adjust the missing-value pattern / values to match the new random split
$ws.Range("F2").Value = ""
$ws.Range("B3").Value = -19.8
$ws.Range("E3").Value = -7.5
$ws.Range("B4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("B6").Value = -18.7
$ws.Range("B7").Value = ""
$ws.Range("F10").Value = 0.7105
$ws.Range("F13").Value = ""
$ws.Range("B15").Value = -18.9
$ws.Range("B16").Value = ""
$ws.Range("F19").Value = 0.71076
$ws.Range("F21").Value = 0.70981
$ws.Range("F22").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("F26").Value = 0.70925
$ws.Range("B27").Value = -19.3
$ws.Range("F27").Value = 0.7092
$ws.Range("B29").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("F30").Value = ""
$ws.Range("B31").Value = -19.5
$ws.Range("B32").Value = ""
$ws.Range("E37").Value = -7.1
$ws.Range("B39").Value = -19.8
$ws.Range("E39").Value = ""
$ws.Range("B40").Value = ""
$ws.Range("F43").Value = 0.71152
$ws.Range("E45").Value = -7.4
$ws.Range("F45").Value = ""
$ws.Range("F46").Value = ""
$ws.Range("E47").Value = ""
$ws.Range("B51").Value = -20.5
$ws.Range("B52").Value = ""
$ws.Range("E53").Value = -5.7
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = 0.71492
$ws.Range("B56").Value = -19.2
$ws.Range("B57").Value = ""
$ws.Range("F58").Value = ""
$ws.Range("E59").Value = -5.7
$ws.Range("E60").Value = -8.1
$ws.Range("F60").Value = 0.70948
$ws.Range("E62").Value = ""
